$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so numeric-looking
# strings (e.g. "1.71", "0.421") are not auto-converted to numbers
# and keep exact formatting (trailing zeros, thousand-dot separators).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "97.682.72"
$ws.Range("E2").Value = "  +2.11%  "
$ws.Range("D3").Value = "3.622.45"
$ws.Range("E3").Value = "  +1.39%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "244.21"
$ws.Range("E5").Value = "  +3.49%  "
$ws.Range("D6").Value = "658.56"
$ws.Range("E6").Value = "  +0.95%  "
$ws.Range("D7").Value = "1.71"
$ws.Range("E7").Value = "  +16.39%  "
$ws.Range("D8").Value = "0.421"
$ws.Range("E8").Value = "  +4.85%  "
$ws.Range("E9").Value = "  +7.23%  "
$ws.Range("D10").Value = "0.999"
$ws.Range("E10").Value = "  -0.03%  "
$ws.Range("D11").Value = "3.621.56"
$ws.Range("E11").Value = "  +1.31%  "
$ws.Range("D12").Value = "43.93"
$ws.Range("E12").Value = "  +3.83%  "
$ws.Range("E13").Value = "  +1.66%  "
$ws.Range("E14").Value = "  -1.75%  "
$ws.Range("D15").Value = "4.292.70"
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("D16").Value = "97.219.43"
$ws.Range("E16").Value = "  +1.79%  "
$ws.Range("E17").Value = "  +2.90%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "8.76"
$ws.Range("E18").Value = "  +10.56%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.619.46"
$ws.Range("E19").Value = "  +1.36%  "
$ws.Range("D20").Value = "12.83"
$ws.Range("E20").Value = "  +0.32%  "
$ws.Range("D21").Value = "18.21"
$ws.Range("E21").Value = "  +1.72%  "
$ws.Range("E22").Value = "  +12.39%  "
$ws.Range("D23").Value = "513.92"
$ws.Range("D24").Value = "3.46"
$ws.Range("E24").Value = "  -1.41%  "
$ws.Range("D25").Value = "0.0000209"
$ws.Range("E25").Value = "  +7.47%  "
$ws.Range("D26").Value = "6.90"
$ws.Range("E26").Value = "  +4.19%  "
$ws.Range("D27").Value = "98.91"
$ws.Range("E27").Value = "  +3.53%  "
$ws.Range("D28").Value = "13.12"
$ws.Range("E28").Value = "  +4.20%  "
$ws.Range("D29").Value = "3.816.76"
$ws.Range("E29").Value = "  +1.45%  "
$ws.Range("E30").Value = "  +9.51%  "
$ws.Range("D31").Value = "3.04"
$ws.Range("E31").Value = "  -0.64%  "
$ws.Range("D32").Value = "11.79"
$ws.Range("E32").Value = "  +4.13%  "
$ws.Range("E33").Value = "  +0.19%  "
$ws.Range("E34").Value = "  +4.20%  "
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("D36").Value = "31.84"
$ws.Range("E36").Value = "  -0.59%  "
$ws.Range("D37").Value = "619.84"
$ws.Range("E37").Value = "  +9.88%  "
$ws.Range("D38").Value = "8.81"
$ws.Range("E38").Value = "  +7.59%  "
$ws.Range("E39").Value = "  +1.89%  "
$ws.Range("E40").Value = "  +8.69%  "
$ws.Range("D41").Value = "1.96"
$ws.Range("E41").Value = "  +10.24%  "
$ws.Range("E42").Value = "  +2.31%  "
$ws.Range("E43").Value = "  +3.29%  "
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("E45").Value = "  +5.02%  "
$ws.Range("E46").Value = "  +5.75%  "
$ws.Range("E47").Value = "  +0.44%  "
$ws.Range("E48").Value = "  +0.42%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "8.56"
$ws.Range("E49").Value = "  +5.92%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "33.40"
$ws.Range("E50").Value = "  -4.86%  "
$ws.Range("D51").Value = "3.55"
$ws.Range("E51").Value = "  -0.84%  "
